# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Alcachofa" (Vega Modelo de Temuco)
# right above the current row 289, pushing the rest of the table down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 289:290 (existing rows 289-329 shift to 291-331)
$ws.Rows("289:290").Insert()

# ---- New row 289: Alcachofa, Española, Primera ----
$ws.Cells.Item(289, 1).Value = 10
$ws.Cells.Item(289, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(289, 3).Value = "La Araucanía"
$ws.Cells.Item(289, 4).Value = 45142
$ws.Cells.Item(289, 5).Value = 9
$ws.Cells.Item(289, 6).Value = 100112013
$ws.Cells.Item(289, 7).Value = "Alcachofa"
$ws.Cells.Item(289, 8).Value = "Española"
$ws.Cells.Item(289, 9).Value = "Primera"
$ws.Cells.Item(289, 10).Value = 50
$ws.Cells.Item(289, 11).Value = 17000
$ws.Cells.Item(289, 12).Value = 17000
$ws.Cells.Item(289, 13).Value = 17000
$ws.Cells.Item(289, 14).Value = "$/caja 35 unidades"
$ws.Cells.Item(289, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(289, 16).Value = 486
$ws.Cells.Item(289, 17).Value = 35
$ws.Cells.Item(289, 18).Value = "Hortaliza"

# ---- New row 290: Alcachofa, Madrigal, Primera ----
$ws.Cells.Item(290, 1).Value = 10
$ws.Cells.Item(290, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(290, 3).Value = "La Araucanía"
$ws.Cells.Item(290, 4).Value = 45142
$ws.Cells.Item(290, 5).Value = 9
$ws.Cells.Item(290, 6).Value = 100112013
$ws.Cells.Item(290, 7).Value = "Alcachofa"
$ws.Cells.Item(290, 8).Value = "Madrigal"
$ws.Cells.Item(290, 9).Value = "Primera"
$ws.Cells.Item(290, 10).Value = 100
$ws.Cells.Item(290, 11).Value = 14000
$ws.Cells.Item(290, 12).Value = 14000
$ws.Cells.Item(290, 13).Value = 14000
$ws.Cells.Item(290, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(290, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(290, 16).Value = 350
$ws.Cells.Item(290, 17).Value = 40
$ws.Cells.Item(290, 18).Value = "Hortaliza"
